# Build site at 2022-09-26 16:07:08 UTC
# Applies the LOT2015.xlsx content restructuring:
#  - the old "Docentes responsaveis" value row (row 13) is removed, which
#    shifts every row below it up by one;
#  - a handful of cells then receive new text content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the standalone row that used to hold the "101761 - Arnaldo Marcio
# Ramalho Prata" value (old row 13); everything below shifts up by one row.
$ws.Rows.Item(13).Delete()

# --- Update the cells whose text changed as part of this edit ---

# Row 10 (Objetivos:) B/C now hold the teacher string instead of the long
# Portuguese objectives paragraph.
$ws.Range("B10").Value = "101761 - Arnaldo Márcio Ramalho Prata"
$ws.Range("C10").Value = "101761 - Arnaldo Márcio Ramalho Prata"

# Row 13 (Programa resumido:) B/C now just say "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (Programa:) B/C now hold the activation date value. Pasted as a
# value (not typed directly) so Excel keeps it as the existing text
# "01/01/2018" instead of auto-converting it to a date serial number.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)

# Row 18 (Método:) B/C now hold the teacher string.
$ws.Range("B18").Value = "101761 - Arnaldo Márcio Ramalho Prata"
$ws.Range("C18").Value = "101761 - Arnaldo Márcio Ramalho Prata"

# Row 19 (Critério:) B/C now hold the previous "Método" evaluation text.
$ws.Range("B19").Value = "Os alunos serão avaliados formalmente por duas provas escritas (P1 e P2), sendo a segunda prova (P2) com peso 2."
$ws.Range("C19").Value = "Os alunos serão avaliados formalmente por duas provas escritas (P1 e P2), sendo a segunda prova (P2) com peso 2."

# Row 20 (Norma de recuperação:) B/C now hold the previous "Critério" text.
$ws.Range("B20").Value = "A nota final (NF) será calculada como: N_F=(P1+(P2×2))/3. Serão aprovados os alunos que obtiverem NF maior ou igual 5,0."
$ws.Range("C20").Value = "A nota final (NF) será calculada como: N_F=(P1+(P2×2))/3. Serão aprovados os alunos que obtiverem NF maior ou igual 5,0."

# Row 21 (Bibliografia:) B/C now hold the previous "Norma de recuperação" text.
$ws.Range("B21").Value = "Será oferecido um programa de recuperação avaliado por uma prova escrita final (PR).
A média de recuperação (MR) será calculada como: MR=(NF+PR)/2. Serão aprovados os alunos que obtiverem MR maior ou igual a 5,0."
$ws.Range("C21").Value = "Será oferecido um programa de recuperação avaliado por uma prova escrita final (PR).
A média de recuperação (MR) será calculada como: MR=(NF+PR)/2. Serão aprovados os alunos que obtiverem MR maior ou igual a 5,0."
